$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83-98 down to 84-99
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new data entry
$ws.Cells.Item(83, 1).Value = 3
$ws.Cells.Item(83, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(83, 3).Value = "Coquimbo"
$ws.Cells.Item(83, 4).Value = 44522
$ws.Cells.Item(83, 5).Value = 5
$ws.Cells.Item(83, 6).Value = 100112026
$ws.Cells.Item(83, 7).Value = "Haba"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 125
$ws.Cells.Item(83, 11).Value = 7000
$ws.Cells.Item(83, 12).Value = 7500
$ws.Cells.Item(83, 13).Value = 7260
$ws.Cells.Item(83, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(83, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(83, 16).Value = 290
$ws.Cells.Item(83, 17).Value = 25
$ws.Cells.Item(83, 18).Value = "Hortaliza"
